$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add "Saldo" header in C1 (plain style, not bold)
$ws.Range("C1").ClearFormats()
$ws.Range("C1").Value = "Saldo"

# Add balance values in C2 and C3
$ws.Range("C2").Value = 1200
$ws.Range("C3").Value = -200

# Reflect the user's selection after adding the new "Saldo" check column
$null = $ws.Range("C1:F1048576").Select()
